$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '22.489.54'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +9.36%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.615.92'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +9.48%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.Style = "Normal"
$ws.Range("E4").Value = '  -0.96%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '304.96'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +9.22%  '

# Row 6
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '0.9904'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.69%  '

# Row 7
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3696'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.24%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '0.3422'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +11.65%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '42.29'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +5.64%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '1.142'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +7.65%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.07096'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +6.34%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.9989'
$c.Style = "Normal"
$ws.Range("E12").Value = '  -0.69%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '19.86'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +9.70%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '5.945'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +7.55%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '6.673'
$c.Style = "Normal"
$ws.Range("E15").Value = '  +7.20%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.00001088'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +5.16%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.605.95'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +8.82%  '

# Row 18
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '0.9905'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +3.72%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06781'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +13.90%  '

# Row 20
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '78.45'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +12.18%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '6.058'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +10.08%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '16.12'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +11.43%  '

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '11.91'
$c.Style = "Normal"
$ws.Range("E23").Value = '  +7.45%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '22.486.48'
$c.Style = "Normal"
$ws.Range("E24").Value = '  +9.02%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.387'
$c.Style = "Normal"
$ws.Range("E25").Value = '  +4.93%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '2.556'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +20.52%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '150.46'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +5.10%  '

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '19.63'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +13.62%  '

# Row 29
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '1.787.92'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +9.26%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '123.10'
$c.Style = "Normal"
$ws.Range("E30").Value = '  +7.91%  '

# Row 31
$ws.Range("B31").Value = 'HuobiToken'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '4.049'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +2.23%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '6.151'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +22.47%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '0.9576'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +17.94%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.08266'
$c.Style = "Normal"
$ws.Range("E34").Value = '  +3.65%  '

# Row 35
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '1.651'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +9.02%  '

# Row 36
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '12.04'
$c.Style = "Normal"
$ws.Range("E36").Value = '  +15.70%  '

# Row 37
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '5.278'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +11.37%  '

# Row 38
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '1.267'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +3.16%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '8.645'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +16.05%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.06112'
$c.Style = "Normal"
$ws.Range("E40").Value = '  +4.51%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.02236'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +8.90%  '

# Row 42
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.2031'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +7.96%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '0.9911'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.70%  '

# Row 44
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.5944'
$c.Style = "Normal"
$ws.Range("E44").Value = '  +11.94%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '3.828'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +8.17%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '13.17'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +6.69%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '0.5726'
$c.Style = "Normal"
$ws.Range("E47").Value = '  +10.10%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '127.26'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +7.83%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '1.988'
$c.Style = "Normal"
$ws.Range("E49").Value = '  +9.14%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '0.06821'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +5.24%  '

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '74.16'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +9.63%  '
